# Move the "Vertical Space" / "Welcome" demo slide (currently slide 24)
# so that it lands right after the slide that currently sits at position 29
# (i.e. it becomes the new slide 29), then retitle it "Headers" since it now
# introduces the upcoming "Headers" section.

$p = $ppt.ActivePresentation

$movedSlide = $p.Slides.Item(24)
$movedSlide.MoveTo(29)

$titleShape = $p.Slides.Item(29).Shapes.Item(1)
$titleShape.TextFrame.TextRange.Text = "Headers"
